$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 1585.909
$ws.Range("I11").Value = 1585.909
$ws.Range("K11").Value = 1585.909
$ws.Range("M11").Value = -1445.909
$ws.Range("H112").Value = 5757.55
$ws.Range("I112").Value = 1442
$ws.Range("J112").Value = 5868.205
$ws.Range("K112").Value = 4326
$ws.Range("L112").Value = 17604.615
$ws.Range("M112").Value = -3218
$ws.Range("N112").Value = -19820.615
$ws.Range("H137").Value = 3111.5881
$ws.Range("I137").Value = 3174.7778
$ws.Range("J137").Value = 3040.5
$ws.Range("K137").Value = 9524.3334
$ws.Range("L137").Value = 9121.5
$ws.Range("M137").Value = -6974.3334
$ws.Range("N137").Value = -14221.5
$ws.Range("H138").Value = 7100.1777
$ws.Range("I138").Value = 2071.5715
$ws.Range("J138").Value = 8026.5
$ws.Range("K138").Value = 6214.7145
$ws.Range("L138").Value = 24079.5
$ws.Range("M138").Value = -1074.7145
$ws.Range("N138").Value = -34359.5
$ws.Range("H141").Value = 1435.3448
$ws.Range("I141").Value = 833.86365
$ws.Range("K141").Value = 2501.59095
$ws.Range("M141").Value = 2678.40905

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1719429.2
$ws.Range("I32").Value = 1817439.5
$ws.Range("K32").Value = 1817439.5
$ws.Range("M32").Value = -1817152.5
$ws.Range("H36").Value = 1000
$ws.Range("I36").Value = 1000
$ws.Range("K36").Value = 1000
$ws.Range("M36").Value = -654
$ws.Range("H45").Value = 5074.391
$ws.Range("I45").Value = 4260.9443
$ws.Range("K45").Value = 4260.9443
$ws.Range("M45").Value = -3883.9443
$ws.Range("H110").Value = 23810906
$ws.Range("I110").Value = 1436.1818
$ws.Range("K110").Value = 1436.1818
$ws.Range("M110").Value = 608.8181999999999
$ws.Range("H122").Value = 3084.7334
$ws.Range("I122").Value = 1781.72
$ws.Range("K122").Value = 5345.16
$ws.Range("M122").Value = -2895.16
$ws.Range("H132").Value = 5191.724
$ws.Range("I132").Value = 4071.606
$ws.Range("K132").Value = 12214.818
$ws.Range("M132").Value = -9684.818000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 100103230
$ws.Range("I86").Value = 333766.66
$ws.Range("J86").Value = 142861570
$ws.Range("K86").Value = 333766.66
$ws.Range("L86").Value = 142861570
$ws.Range("M86").Value = -332643.66
$ws.Range("N86").Value = -142863816
$ws.Range("H89").Value = 100103230
$ws.Range("I89").Value = 333766.66
$ws.Range("J89").Value = 142861570
$ws.Range("K89").Value = 1668833.3
$ws.Range("L89").Value = 714307850
$ws.Range("M89").Value = -1663217.3
$ws.Range("N89").Value = -714319082

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5669.589
$ws.Range("I31").Value = 2621.1428
$ws.Range("K31").Value = 2621.1428
$ws.Range("M31").Value = -2326.1428
$ws.Range("H34").Value = 5669.589
$ws.Range("I34").Value = 2621.1428
$ws.Range("K34").Value = 2621.1428
$ws.Range("M34").Value = -2419.1428
$ws.Range("H105").Value = 7940852
$ws.Range("I105").Value = 11906445
$ws.Range("K105").Value = 11906445
$ws.Range("M105").Value = -11904698
$ws.Range("H107").Value = 2178.7222
$ws.Range("I107").Value = 567.375
$ws.Range("K107").Value = 567.375
$ws.Range("M107").Value = 1352.625
$ws.Range("H134").Value = 4364.2793
$ws.Range("I134").Value = 1953.4286
$ws.Range("K134").Value = 5860.2858
$ws.Range("M134").Value = -3325.2858

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 16667453
$ws.Range("I14").Value = 16667453
$ws.Range("K14").Value = 50002359
$ws.Range("M14").Value = -50002186
$ws.Range("H134").Value = 76836.21
$ws.Range("I134").Value = 87975.586
$ws.Range("K134").Value = 263926.758
$ws.Range("M134").Value = -258856.758

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H52").Value = 64999.6
$ws.Range("I52").Value = 44999
$ws.Range("K52").Value = 44999
$ws.Range("M52").Value = -44740
$ws.Range("H80").Value = 1948.8462
$ws.Range("I80").Value = 1794.2
$ws.Range("J80").Value = 2464.3333
$ws.Range("K80").Value = 1794.2
$ws.Range("L80").Value = 2464.3333
$ws.Range("M80").Value = -796.2
$ws.Range("N80").Value = -4460.3333
$ws.Range("H83").Value = 1948.8462
$ws.Range("I83").Value = 1794.2
$ws.Range("J83").Value = 2464.3333
$ws.Range("K83").Value = 8971
$ws.Range("L83").Value = 12321.6665
$ws.Range("M83").Value = -3979
$ws.Range("N83").Value = -22305.6665
$ws.Range("H113").Value = 300177.25
$ws.Range("I113").Value = 771711.7
$ws.Range("J113").Value = 8275
$ws.Range("K113").Value = 771711.7
$ws.Range("L113").Value = 8275
$ws.Range("M113").Value = -769541.7
$ws.Range("N113").Value = -12615

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2300.7856
$ws.Range("I22").Value = 1012.5
$ws.Range("J22").Value = 4018.5
$ws.Range("K22").Value = 1012.5
$ws.Range("L22").Value = 4018.5
$ws.Range("M22").Value = -717.5
$ws.Range("N22").Value = -4608.5
$ws.Range("H27").Value = 2300.7856
$ws.Range("I27").Value = 1012.5
$ws.Range("J27").Value = 4018.5
$ws.Range("K27").Value = 1012.5
$ws.Range("L27").Value = 4018.5
$ws.Range("M27").Value = -905.5
$ws.Range("N27").Value = -4232.5
$ws.Range("H31").Value = 2108
$ws.Range("I31").Value = 628.25
$ws.Range("J31").Value = 4081
$ws.Range("K31").Value = 628.25
$ws.Range("L31").Value = 4081
$ws.Range("M31").Value = -380.25
$ws.Range("N31").Value = -4577
$ws.Range("H46").Value = 5295040
$ws.Range("J46").Value = 5852149
$ws.Range("L46").Value = 5852149
$ws.Range("N46").Value = -5852525
$ws.Range("H55").Value = 378.83334
$ws.Range("I55").Value = 106.111115
$ws.Range("K55").Value = 106.111115
$ws.Range("M55").Value = 66.888885
$ws.Range("H136").Value = 7486.278
$ws.Range("I136").Value = 1529.0714
$ws.Range("J136").Value = 11277.228
$ws.Range("K136").Value = 4587.2142
$ws.Range("L136").Value = 33831.68399999999
$ws.Range("M136").Value = -2037.2142
$ws.Range("N136").Value = -38931.68399999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 22223562
$ws.Range("I107").Value = 1357.8572
$ws.Range("J107").Value = 41667988
$ws.Range("K107").Value = 4073.5716
$ws.Range("L107").Value = 125003964
$ws.Range("M107").Value = -2153.5716
$ws.Range("N107").Value = -125007804
$ws.Range("H122").Value = 204312.34
$ws.Range("I122").Value = 402389.4
$ws.Range("J122").Value = 6235.3
$ws.Range("K122").Value = 1207168.2
$ws.Range("L122").Value = 18705.9
$ws.Range("M122").Value = -1204718.2
$ws.Range("N122").Value = -23605.9
$ws.Range("H123").Value = 50399.5
$ws.Range("J123").Value = 50399.5
$ws.Range("L123").Value = 50399.5
$ws.Range("N123").Value = -60199.5
